# Add option to one hot encode zoom levels
# - Fill in the "File" column for the existing row 54 (previously blank)
# - Add a new table row (row 61) describing the "incl. One hot Zoom levels" run
# - Keep conditional formatting ranges and table range in sync with the new row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item(1)

# Grow the table by one row (this expands table ref, autoFilter ref and the
# worksheet dimension automatically)
$tbl.ListRows.Add() | Out-Null

# Set the "Comments" (I61) and "Average Val AUC" (M61) cells first, then the
# "File" cells for row 54 and row 61, so that new shared strings are created
# in the same order as the target workbook.
$ws.Range("I61").Value2 = "incl. One hot Zoom levels"
$ws.Range("M61").Value2 = "0.667 (0.053)"

$ws.Range("A54").Value2 = "2023-03-12-1553_RF_zoom_avg.csv"
$ws.Range("A61").Value2 = "2023-03-12-1624_RF_zoomonehot_avg.csv"

$ws.Range("B61").Value2 = "RandomForest"
$ws.Range("C61").Value2 = "MoCo"
$ws.Range("D61").Value2 = "Centers"
$ws.Range("E61").Value2 = "1 x 3"
$ws.Range("F61").Value2 = "average"
$ws.Range("J61").Value2 = 0.623
$ws.Range("K61").Value2 = 0.742
$ws.Range("L61").Value2 = 0.637

# Extend the top-10 conditional formatting on columns J and K so it still
# covers the whole data range (previously J2:J60 / K2:K60).
$oldK = '$K$2:$K$60'
$oldJ = '$J$2:$J$60'
$allConditions = $ws.Cells.FormatConditions
for ($i = 1; $i -le $allConditions.Count; $i++) {
    $cond = $allConditions.Item($i)
    $addr = $cond.AppliesTo.Address()
    if ($addr -eq $oldK) {
        $cond.ModifyAppliesToRange($ws.Range("K2:K61"))
    }
    elseif ($addr -eq $oldJ) {
        $cond.ModifyAppliesToRange($ws.Range("J2:J61"))
    }
}

# Match the updated view state from the authored workbook: scrolled up to
# show the new row, with cell A62 selected.
$ws.Range("A62").Select() | Out-Null
